$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing header cell (H1) onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns (I0 / IF) for rows 2-8
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 9

$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 9
